$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "total_calories_burned"
$ws.Range("C1").Value = "daily_step_count"

$ws.Range("A1").Select()
